$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A to hold the new "model_id" field
$ws.Range("A1:A7").EntireColumn.Insert()

# Header for new column
$ws.Range("A1").Value = "model_id"

# model_id values
$ws.Range("A2").Value = 10
$ws.Range("A3").Value = 20
$ws.Range("A4").Value = 30
$ws.Range("A5").Value = 40
$ws.Range("A6").Value = 50
$ws.Range("A7").Value = 60

# Update selection to match recorded state
$ws.Range("A8").Select()
